$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-07-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-24 Monday", 2) | Out-Null

# Update each answer cell in the multiplication table (row-major order)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "85×36=3060"
$t.Cell(1, 2).Range.Text = "89×41=3649"
$t.Cell(1, 3).Range.Text = "49×99=4851"
$t.Cell(1, 4).Range.Text = "55×82=4510"
$t.Cell(1, 5).Range.Text = "57×45=2565"
$t.Cell(2, 1).Range.Text = "21×66=1386"
$t.Cell(2, 2).Range.Text = "99×70=6930"
$t.Cell(2, 3).Range.Text = "54×71=3834"
$t.Cell(2, 4).Range.Text = "58×18=1044"
$t.Cell(2, 5).Range.Text = "72×19=1368"
$t.Cell(3, 1).Range.Text = "74×55=4070"
$t.Cell(3, 2).Range.Text = "10×57=570"
$t.Cell(3, 3).Range.Text = "33×67=2211"
$t.Cell(3, 4).Range.Text = "67×50=3350"
$t.Cell(3, 5).Range.Text = "60×41=2460"
$t.Cell(4, 1).Range.Text = "61×60=3660"
$t.Cell(4, 2).Range.Text = "91×98=8918"
$t.Cell(4, 3).Range.Text = "70×92=6440"
$t.Cell(4, 4).Range.Text = "15×86=1290"
$t.Cell(4, 5).Range.Text = "78×71=5538"
$t.Cell(5, 1).Range.Text = "11×43=473"
$t.Cell(5, 2).Range.Text = "24×69=1656"
$t.Cell(5, 3).Range.Text = "89×50=4450"
$t.Cell(5, 4).Range.Text = "83×35=2905"
$t.Cell(5, 5).Range.Text = "57×61=3477"
$t.Cell(6, 1).Range.Text = "95×50=4750"
$t.Cell(6, 2).Range.Text = "40×49=1960"
$t.Cell(6, 3).Range.Text = "61×41=2501"
$t.Cell(6, 4).Range.Text = "53×38=2014"
$t.Cell(6, 5).Range.Text = "37×14=518"
$t.Cell(7, 1).Range.Text = "38×48=1824"
$t.Cell(7, 2).Range.Text = "33×47=1551"
$t.Cell(7, 3).Range.Text = "32×43=1376"
$t.Cell(7, 4).Range.Text = "19×77=1463"
$t.Cell(7, 5).Range.Text = "77×21=1617"
$t.Cell(8, 1).Range.Text = "74×24=1776"
$t.Cell(8, 2).Range.Text = "49×20=980"
$t.Cell(8, 3).Range.Text = "22×31=682"
$t.Cell(8, 4).Range.Text = "47×97=4559"
$t.Cell(8, 5).Range.Text = "48×24=1152"
$t.Cell(9, 1).Range.Text = "95×96=9120"
$t.Cell(9, 2).Range.Text = "24×24=576"
$t.Cell(9, 3).Range.Text = "21×24=504"
$t.Cell(9, 4).Range.Text = "19×29=551"
$t.Cell(9, 5).Range.Text = "16×97=1552"
$t.Cell(10, 1).Range.Text = "65×29=1885"
$t.Cell(10, 2).Range.Text = "54×56=3024"
$t.Cell(10, 3).Range.Text = "73×83=6059"
$t.Cell(10, 4).Range.Text = "35×85=2975"
$t.Cell(10, 5).Range.Text = "28×100=2800"
$t.Cell(11, 1).Range.Text = "25×97=2425"
$t.Cell(11, 2).Range.Text = "45×39=1755"
$t.Cell(11, 3).Range.Text = "83×65=5395"
$t.Cell(11, 4).Range.Text = "85×22=1870"
$t.Cell(11, 5).Range.Text = "73×40=2920"
$t.Cell(12, 1).Range.Text = "29×32=928"
$t.Cell(12, 2).Range.Text = "87×57=4959"
$t.Cell(12, 3).Range.Text = "99×95=9405"
$t.Cell(12, 4).Range.Text = "75×45=3375"
$t.Cell(12, 5).Range.Text = "39×30=1170"
$t.Cell(13, 1).Range.Text = "35×57=1995"
$t.Cell(13, 2).Range.Text = "97×52=5044"
$t.Cell(13, 3).Range.Text = "39×76=2964"
$t.Cell(13, 4).Range.Text = "74×63=4662"
$t.Cell(13, 5).Range.Text = "46×100=4600"
$t.Cell(14, 1).Range.Text = "76×19=1444"
$t.Cell(14, 2).Range.Text = "78×42=3276"
$t.Cell(14, 3).Range.Text = "48×55=2640"
$t.Cell(14, 4).Range.Text = "16×68=1088"
$t.Cell(14, 5).Range.Text = "73×95=6935"
$t.Cell(15, 1).Range.Text = "55×24=1320"
$t.Cell(15, 2).Range.Text = "76×95=7220"
$t.Cell(15, 3).Range.Text = "18×70=1260"
$t.Cell(15, 4).Range.Text = "79×20=1580"
$t.Cell(15, 5).Range.Text = "51×61=3111"
$t.Cell(16, 1).Range.Text = "83×11=913"
$t.Cell(16, 2).Range.Text = "59×75=4425"
$t.Cell(16, 3).Range.Text = "58×56=3248"
$t.Cell(16, 4).Range.Text = "84×99=8316"
$t.Cell(16, 5).Range.Text = "57×77=4389"
$t.Cell(17, 1).Range.Text = "32×53=1696"
$t.Cell(17, 2).Range.Text = "56×45=2520"
$t.Cell(17, 3).Range.Text = "48×98=4704"
$t.Cell(17, 4).Range.Text = "94×28=2632"
$t.Cell(17, 5).Range.Text = "66×70=4620"
$t.Cell(18, 1).Range.Text = "19×97=1843"
$t.Cell(18, 2).Range.Text = "25×40=1000"
$t.Cell(18, 3).Range.Text = "61×99=6039"
$t.Cell(18, 4).Range.Text = "38×95=3610"
$t.Cell(18, 5).Range.Text = "54×32=1728"
$t.Cell(19, 1).Range.Text = "79×64=5056"
$t.Cell(19, 2).Range.Text = "58×24=1392"
$t.Cell(19, 3).Range.Text = "31×39=1209"
$t.Cell(19, 4).Range.Text = "33×91=3003"
$t.Cell(19, 5).Range.Text = "82×84=6888"
$t.Cell(20, 1).Range.Text = "28×80=2240"
$t.Cell(20, 2).Range.Text = "88×12=1056"
$t.Cell(20, 3).Range.Text = "34×12=408"
$t.Cell(20, 4).Range.Text = "71×70=4970"
$t.Cell(20, 5).Range.Text = "28×44=1232"
